$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A49").Value = 0
$ws.Range("B49").Value = 0
$ws.Range("C49").Value = 0
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0.090909
$ws.Range("F49").Value = -2.203442372573805
$ws.Range("G49").Value = "query"

$ws.Range("A50").Value = 0
$ws.Range("B50").Value = 0
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0.272727
$ws.Range("F50").Value = -2.382709596165335
$ws.Range("G50").Value = "query"
